$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> columns that flip from 0 to 1
$updates = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("D", "E")
    6  = @("H")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("D", "E")
    11 = @("H")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("D", "E")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
